$d = $word.ActiveDocument

# Append a new numbered-list paragraph after the last one, inheriting its
# style/numbering (pStyle "a3" + numPr ilvl 0 / numId 1).
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$lastPara.Range.InsertParagraphAfter()

$newIndex = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($newIndex)

# Temporary trailing sentinel character lets us collapse a Range to a point
# that sits *inside* a run (not exactly on a run boundary) before adding the
# _GoBack bookmark, then we trim it back off.
$newPara.Range.Text = "Антон, проверь, плез, на наличие ненужных файловX"
$newPara.Range.LanguageID = "en-US"

$tail = $newPara.Range
$null = $tail.MoveEnd(1, -2)
$tail.Collapse(0)
$d.Bookmarks.Add("_GoBack", $tail)

$null = $tail.MoveEnd(1, 1)
$tail.Text = ""
